$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (date 44943) is inserted before the existing
# row 380, pushing that row and everything below it down by one row
# (old row 380 -> new row 381, ..., old row 463 -> new row 464).
$ws.Rows.Item(380).Insert()

$ws.Range("A380").Value2 = 8
$ws.Range("B380").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C380").Value2 = "Coquimbo"
$ws.Range("D380").Value2 = 44943
$ws.Range("E380").Value2 = 4
$ws.Range("F380").Value2 = 100114013
$ws.Range("G380").Value2 = "Zanahoria"
$ws.Range("H380").Value2 = "Sin especificar"
$ws.Range("I380").Value2 = "Primera"
$ws.Range("J380").Value2 = 400
$ws.Range("K380").Value2 = 6000
$ws.Range("L380").Value2 = 6500
$ws.Range("M380").Value2 = 6250
$ws.Range("N380").Value2 = "`$/saco 20 kilos"
$ws.Range("O380").Value2 = "Provincia del Elquí"
$ws.Range("P380").Value2 = 312
$ws.Range("Q380").Value2 = 20
$ws.Range("R380").Value2 = "Hortaliza"
